$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Make room for 5 new "Import" rows (Spring expression classes) right
#    after the existing imports / before the old "Functions" row.
# ---------------------------------------------------------------------------
$ws.Range("A8:A12").EntireRow.Insert()

# Copy the formatting of the last existing Import row (row 6) onto the 5
# freshly inserted rows (7-11) so they look like the other Import rows.
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D11").PasteSpecial(-4122)
$ws.Range("A7:D11").RowHeight = 13.3

# ---------------------------------------------------------------------------
# 2. New Import rows: Spring expression / SpEL classes
# ---------------------------------------------------------------------------
$ws.Range("C7").Value = "Import"
$ws.Range("D7").Value = "org.springframework.expression.EvaluationContext"

$ws.Range("C8").Value = "Import"
$ws.Range("D8").Value = "org.springframework.expression.Expression"

$ws.Range("C9").Value = "Import"
$ws.Range("D9").Value = "org.springframework.expression.ExpressionParser"

$ws.Range("C10").Value = "Import"
$ws.Range("D10").Value = "org.springframework.expression.spel.standard.SpelExpressionParser"

$ws.Range("C11").Value = "Import"
$ws.Range("D11").Value = "org.springframework.expression.spel.support.StandardEvaluationContext"

# ---------------------------------------------------------------------------
# 3. Row 12 ("Functions") - append a new evalSpring() helper function to the
#    existing dateFormat() one, as rich text (two runs: plain + styled).
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = "Functions"

$funcText1 = @'
function String dateFormat(String fmt)
{
  return LocalDate.now().toString(DateTimeFormat.forPattern(fmt));
}

'@

$funcText2 = @'

function Boolean evalSpring(String expression, Object obj)
{
    ExpressionParser ep = new SpelExpressionParser();
    Expression exp = ep.parseExpression(expression);
    EvaluationContext ec = new StandardEvaluationContext();
    Boolean evaluated = exp.getValue(ec, obj, Boolean.class);
    return evaluated;
}
'@

$ws.Range("D12").Value = $funcText1 + $funcText2
$ws.Range("D12").RowHeight = 169.4

$run2Start = $funcText1.Length + 1
$run2Len = $funcText2.Length
$run2 = $ws.Range("D12").Characters($run2Start, $run2Len)
$run2.Font.Name = "Calibri"
$run2.Font.Size = 11
$run2.Font.Color = 0

# ---------------------------------------------------------------------------
# 4. Rule-table condition / rule-name text updates.
# ---------------------------------------------------------------------------

# Condition expression now delegates to the new evalSpring() helper.
$ws.Range("C18").Value = 'eval(evalSpring("$param", $complaint))'
$ws.Range("C18").Font.Name = "Calibri"
$ws.Range("C18").Font.Size = 11
$ws.Range("C18").Font.Color = 0

$ws.Range("C18").RowHeight = 13.3

# "Field is null" -> "expression is true"
$ws.Range("C19").Value = "expression is true"

# "complaintNumber" -> "complaintNumber == null"
$ws.Range("C20").Value = "complaintNumber == null"

# "containerFolder.cmisFolderId" -> "container?.folder?.cmisFolderId == null"
$ws.Range("C21").Value = "container?.folder?.cmisFolderId == null"

# ---------------------------------------------------------------------------
# 5. Selection / view bookkeeping to match the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("C20").Select()
